$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of spreadsheet row number -> corrected sex value for column A
# (re-ran anonymisation fix: swap mis-mapped Male/Female codes for these rows)
$rowValues = @{
    4 = "Female"
    10 = "Male"
    11 = "Female"
    14 = "Female"
    15 = "Female"
    25 = "Female"
    27 = "Female"
    30 = "Male"
    37 = "Male"
    38 = "Male"
    44 = "Male"
    46 = "Female"
    49 = "Female"
    58 = "Male"
    62 = "Male"
    63 = "Female"
    69 = "Female"
    71 = "Male"
    78 = "Female"
    80 = "Male"
    82 = "Male"
    84 = "Male"
    92 = "Male"
    103 = "Female"
    105 = "Female"
    107 = "Female"
    111 = "Female"
    115 = "Male"
    116 = "Male"
    117 = "Male"
    118 = "Female"
    121 = "Male"
    124 = "Female"
    125 = "Male"
    127 = "Male"
    129 = "Female"
    139 = "Male"
    145 = "Male"
    148 = "Female"
    161 = "Female"
    162 = "Male"
    165 = "Male"
    166 = "Male"
    174 = "Female"
    178 = "Female"
    181 = "Female"
    182 = "Female"
    183 = "Male"
    184 = "Male"
    185 = "Male"
    189 = "Female"
    190 = "Male"
    192 = "Female"
    195 = "Female"
    198 = "Male"
    200 = "Female"
}

foreach ($row in $rowValues.Keys) {
    $ws.Cells.Item($row, 1).Value = $rowValues[$row]
}
